$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# id_corso (column H) values for rows 5-14
$ws.Range("H5").Value = 1
$ws.Range("H6").Value = 1
$ws.Range("H7").Value = 2
$ws.Range("H8").Value = 2
$ws.Range("H9").Value = 2
$ws.Range("H10").Value = 3
$ws.Range("H11").Value = 3
$ws.Range("H12").Value = 3
$ws.Range("H13").Value = 3
$ws.Range("H14").Value = 1

# docente (column I) values newly populated
$ws.Range("I7").Value = "petrignani"
$ws.Range("I14").Value = "petrignani"

# I17 becomes a present-but-empty, underlined cell (matches the
# "docente not yet assigned" placeholder style used elsewhere in the
# sheet, e.g. I4), without putting any value into it.
$ws.Range("I17").Font.Underline = 2
$ws.Range("I17").Font.Name = "Arial"
$ws.Range("I17").Font.Size = 10
$ws.Range("I17").Font.Color = $ws.Range("I4").Font.Color

# Move the active selection to I14, matching the author's final cursor position
$ws.Range("I14").Select()
